# Apply the changes described by the commit:
#   "Userdata gets uploaded in DB, Events get created, Tokens are generated"
#
# The worksheet lists a contact (Rafael Fiedler) together with a generated
# token/handle e-mail address in column D (a hyperlinked mailto: link).
# The generated address changed from 6136@htl.rennweg.at to
# 6138@htl.rennweg.at. Updating the cell value automatically reorders the
# shared-string table the same way Excel does (the old, now-unused string
# is dropped and the new one is appended), which also shifts the shared
# string indices referenced by B1/C1/D1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "6138@htl.rennweg.at"

# The active selection moved from C2 to D2 when the workbook was last saved.
$ws.Range("D2").Select() | Out-Null
